# Adds a new "2023" column (T) to the statistics table, mirroring the
# existing "2022" column (S) for number format / font / border, and fills
# in the 2023 data values for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> 2023 value (column T), mirroring column S (2022) row-for-row.
$newValues = @{
    4  = 2023
    5  = 4.8187602774004432
    6  = 11.788953009068425
    7  = 5.2855407047387608
    8  = 11.35112240576027
    9  = 16.577540106951872
    10 = 14.651002073255009
    11 = 5.034965034965035
    12 = 3.1837160751565765
    13 = 2.2263731825525039
    14 = 5.1321450522433931
}

foreach ($row in 4..14) {
    # Column S = 19, column T = 20. Copy S's formatting (number format,
    # font, border) onto T first, then overwrite with the 2023 value.
    $src = $ws.Cells.Item($row, 19)
    $dst = $ws.Cells.Item($row, 20)
    $src.Copy($dst)
    $dst.Value = $newValues[$row]
}
